$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.213.75"
$ws.Range("E2").Value = "  -8.03%  "

$ws.Range("D3").Value = "2.877.40"
$ws.Range("E3").Value = "  -10.89%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'475.16"
$ws.Range("E5").Value = "  -12.11%  "

$ws.Range("D6").Value = "'126.09"
$ws.Range("E6").Value = "  -7.42%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "2.875.86"
$ws.Range("E8").Value = "  -10.91%  "

$ws.Range("D9").Value = "'0.404"
$ws.Range("E9").Value = "  -12.16%  "

$ws.Range("E10").Value = "  -12.39%  "

$ws.Range("D11").Value = "'0.0975"
$ws.Range("E11").Value = "  -15.13%  "

$ws.Range("E12").Value = "  -15.34%  "

$ws.Range("E13").Value = "  -3.89%  "

$ws.Range("D14").Value = "3.372.47"
$ws.Range("E14").Value = "  -10.85%  "

$ws.Range("D15").Value = "'22.82"
$ws.Range("E15").Value = "  -12.19%  "

$ws.Range("D16").Value = "54.198.57"
$ws.Range("E16").Value = "  -8.11%  "

$ws.Range("D17").Value = "2.884.03"
$ws.Range("E17").Value = "  -10.71%  "

$ws.Range("E18").Value = "  -14.67%  "

$ws.Range("D19").Value = "'5.25"
$ws.Range("E19").Value = "  -11.09%  "

$ws.Range("D20").Value = "'11.61"
$ws.Range("E20").Value = "  -13.20%  "

$ws.Range("D21").Value = "'7.11"
$ws.Range("E21").Value = "  -13.65%  "

$ws.Range("D22").Value = "'309.50"
$ws.Range("E22").Value = "  -14.73%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").Value = "'0.449"
$ws.Range("E24").Value = "  -13.99%  "

$ws.Range("D25").Value = "'59.73"
$ws.Range("E25").Value = "  -15.37%  "

$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("D27").Value = "'0.153"
$ws.Range("E27").Value = "  -10.33%  "

$ws.Range("D29").Value = "0.0₃0823"
$ws.Range("E29").Value = "  -15.30%  "

$ws.Range("E30").Value = "  -12.27%  "

$ws.Range("E31").Value = "  -5.90%  "

$ws.Range("D32").Value = "'6.23"
$ws.Range("E32").Value = "  -12.49%  "

$ws.Range("D33").Value = "'19.12"
$ws.Range("E33").Value = "  -12.89%  "

$ws.Range("D34").Value = "'1.62"
$ws.Range("E34").Value = "  -16.23%  "

$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  -14.01%  "

$ws.Range("D36").Value = "'139.35"
$ws.Range("E36").Value = "  -14.13%  "

$ws.Range("D37").Value = "'5.47"
$ws.Range("E37").Value = "  -14.99%  "

$ws.Range("E38").Value = "  -15.85%  "

$ws.Range("D39").Value = "'23.01"
$ws.Range("E39").Value = "  -12.83%  "

$ws.Range("E40").Value = "  -12.37%  "

$ws.Range("D41").Value = "2.904.80"
$ws.Range("E41").Value = "  -10.85%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "'35.41"
$ws.Range("E43").Value = "  -13.88%  "

$ws.Range("D44").Value = "'0.962"
$ws.Range("E44").Value = "  -12.89%  "

$ws.Range("D45").Value = "'0.601"
$ws.Range("E45").Value = "  -16.10%  "

$ws.Range("D46").Value = "'3.43"
$ws.Range("E46").Value = "  -15.06%  "

$ws.Range("D47").Value = "'1.32"
$ws.Range("E47").Value = "  -12.61%  "

$ws.Range("D48").Value = "2.060.34"
$ws.Range("E48").Value = "  -10.49%  "

$ws.Range("D49").Value = "'5.35"
$ws.Range("E49").Value = "  -15.51%  "

$ws.Range("D50").Value = "'17.97"
$ws.Range("E50").Value = "  -14.21%  "

$ws.Range("E51").Value = "  -11.91%  "
